$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.931.97'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '3.409.43'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''408.19'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '''128.34'
$ws.Range("E6").Value = '  -4.39%  '
$ws.Range("D7").Value = '''0.632'
$ws.Range("E7").Value = '  +6.58%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.732'
$ws.Range("E9").Value = '  +6.64%  '
$ws.Range("E10").Value = '  +17.00%  '
$ws.Range("D11").Value = '''42.29'
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").Value = '''0.0000218'
$ws.Range("E12").Value = '  +68.17%  '
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '3.955.22'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '''8.91'
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("D16").Value = '''20.70'
$ws.Range("E16").Value = '  +4.19%  '
$ws.Range("D17").Value = '3.423.08'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '''12.07'
$ws.Range("E18").Value = '  +9.55%  '
$ws.Range("D19").Value = '''1.07'
$ws.Range("E19").Value = '  +5.32%  '
$ws.Range("D20").Value = '61.872.89'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '''413.61'
$ws.Range("E21").Value = '  +31.44%  '
$ws.Range("D22").Value = '''88.99'
$ws.Range("E22").Value = '  +5.72%  '
$ws.Range("D23").Value = '''3.16'
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").Value = '''33.01'
$ws.Range("E26").Value = '  +11.61%  '
$ws.Range("E27").Value = '  +7.94%  '
$ws.Range("D28").Value = '''4.77'
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").Value = '''7.57'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '''2.68'
$ws.Range("E30").Value = '  -4.65%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '''11.86'
$ws.Range("E31").Value = '  +4.21%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '''0.171'
$ws.Range("E32").Value = '  -1.99%  '
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").Value = '''42.66'
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = '''0.0495'
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("D37").Value = '''54.06'
$ws.Range("E37").Value = '  +4.56%  '
$ws.Range("D38").Value = '''0.998'
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("E40").Value = '  +6.23%  '
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("E42").Value = '  +3.91%  '
$ws.Range("D43").Value = '''141.41'
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("E46").Value = '  +8.48%  '
$ws.Range("D47").Value = '''16.56'
$ws.Range("E47").Value = '  -0.95%  '
$ws.Range("D48").Value = '''21.72'
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("D49").Value = '2.107.75'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("E50").Value = '  +2.26%  '
$ws.Range("D51").Value = '''0.130'
$ws.Range("E51").Value = '  +14.19%  '
